# NYPD CompStat weekly report update
# - Bumps the report "Volume/Number" header and the covered-week date range
#   (both are rich-text shared strings; only specific runs change).
# - Refreshes the crime-complaint figures (rows 14-30, columns C:N) with the
#   newly collected week's numbers.
# - Row 30 ("Hate Crimes") additionally changes which columns hold numeric
#   data vs. the "0"/"***.*" placeholder text, so those three cells need an
#   explicit NumberFormat swap alongside the value change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rich-text runs -------------------------------------------------
# "Volume 30   Number  2" -> "...Number  3" (only the trailing run changes)
$ws.Range("A8").Characters(21, 1).Text = "3"

# "Report Covering the Week  1/9/2023  Through  1/15/2023"
#                         -> "...1/16/2023  Through  1/22/2023"
# Edit the right-hand run first so the left-hand run's character offset
# (computed against the original string) is not shifted by a length change.
$ws.Range("C9").Characters(46, 9).Text = "1/22/2023"
$ws.Range("C9").Characters(27, 8).Text = "1/16/2023"

# --- Weekly crime-complaint figures (rows 14-30) ---------------------------
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 5
$ws.Range("H14").Value2 = -80
$ws.Range("J14").Value2 = 3
$ws.Range("C15").Value2 = 2
$ws.Range("D15").Value2 = 2
$ws.Range("E15").Value2 = 0
$ws.Range("F15").Value2 = 14
$ws.Range("G15").Value2 = 10
$ws.Range("H15").Value2 = 40
$ws.Range("I15").Value2 = 11
$ws.Range("J15").Value2 = 10
$ws.Range("K15").Value2 = 10
$ws.Range("L15").Value2 = 57.142857142857
$ws.Range("M15").Value2 = 175
$ws.Range("N15").Value2 = -38.888888888888
$ws.Range("C16").Value2 = 15
$ws.Range("D16").Value2 = 29
$ws.Range("E16").Value2 = -48.275862068965
$ws.Range("F16").Value2 = 82
$ws.Range("G16").Value2 = 115
$ws.Range("H16").Value2 = -28.695652173913
$ws.Range("I16").Value2 = 59
$ws.Range("J16").Value2 = 90
$ws.Range("K16").Value2 = -34.444444444444
$ws.Range("L16").Value2 = -7.8125
$ws.Range("M16").Value2 = -65.088757396449
$ws.Range("N16").Value2 = -88.909774436090
$ws.Range("C17").Value2 = 52
$ws.Range("D17").Value2 = 41
$ws.Range("E17").Value2 = 26.829268292682
$ws.Range("F17").Value2 = 198
$ws.Range("G17").Value2 = 170
$ws.Range("H17").Value2 = 16.470588235294
$ws.Range("I17").Value2 = 157
$ws.Range("J17").Value2 = 136
$ws.Range("K17").Value2 = 15.441176470588
$ws.Range("L17").Value2 = 20.769230769230
$ws.Range("M17").Value2 = 96.25
$ws.Range("N17").Value2 = -25.943396226415
$ws.Range("C18").Value2 = 14
$ws.Range("D18").Value2 = 22
$ws.Range("E18").Value2 = -36.363636363636
$ws.Range("F18").Value2 = 79
$ws.Range("G18").Value2 = 96
$ws.Range("H18").Value2 = -17.708333333333
$ws.Range("I18").Value2 = 68
$ws.Range("J18").Value2 = 70
$ws.Range("K18").Value2 = -2.857142857142
$ws.Range("L18").Value2 = -2.857142857142
$ws.Range("M18").Value2 = -45.6
$ws.Range("N18").Value2 = -88.435374149659
$ws.Range("C19").Value2 = 57
$ws.Range("D19").Value2 = 72
$ws.Range("E19").Value2 = -20.833333333333
$ws.Range("F19").Value2 = 247
$ws.Range("G19").Value2 = 309
$ws.Range("H19").Value2 = -20.064724919093
$ws.Range("I19").Value2 = 195
$ws.Range("J19").Value2 = 232
$ws.Range("K19").Value2 = -15.948275862069
$ws.Range("L19").Value2 = 62.5
$ws.Range("M19").Value2 = 24.203821656051
$ws.Range("N19").Value2 = -42.647058823529
$ws.Range("C20").Value2 = 24
$ws.Range("D20").Value2 = 41
$ws.Range("E20").Value2 = -41.463414634146
$ws.Range("F20").Value2 = 139
$ws.Range("G20").Value2 = 127
$ws.Range("H20").Value2 = 9.448818897637
$ws.Range("I20").Value2 = 102
$ws.Range("J20").Value2 = 100
$ws.Range("K20").Value2 = 2
$ws.Range("L20").Value2 = 78.947368421052
$ws.Range("M20").Value2 = -11.304347826087
$ws.Range("N20").Value2 = -91.528239202657
$ws.Range("C21").Value2 = 164
$ws.Range("E21").Value2 = -21.153846153846
$ws.Range("F21").Value2 = 760
$ws.Range("G21").Value2 = 832
$ws.Range("H21").Value2 = -8.653846153846
$ws.Range("I21").Value2 = 592
$ws.Range("J21").Value2 = 641
$ws.Range("K21").Value2 = -7.644305772230
$ws.Range("L21").Value2 = 31.555555555555
$ws.Range("M21").Value2 = -9.480122324159
$ws.Range("N21").Value2 = -79.600275671950
$ws.Range("C22").Value2 = 4
$ws.Range("D22").Value2 = 4
$ws.Range("E22").Value2 = 0
$ws.Range("F22").Value2 = 12
$ws.Range("G22").Value2 = 11
$ws.Range("H22").Value2 = 9.090909090909
$ws.Range("I22").Value2 = 8
$ws.Range("J22").Value2 = 9
$ws.Range("K22").Value2 = -11.111111111111
$ws.Range("L22").Value2 = 33.333333333333
$ws.Range("M22").Value2 = 33.333333333333
$ws.Range("C23").Value2 = 2
$ws.Range("D23").Value2 = 2
$ws.Range("E23").Value2 = 0
$ws.Range("G23").Value2 = 18
$ws.Range("H23").Value2 = -22.222222222222
$ws.Range("I23").Value2 = 13
$ws.Range("J23").Value2 = 15
$ws.Range("K23").Value2 = -13.333333333333
$ws.Range("L23").Value2 = 85.714285714285
$ws.Range("M23").Value2 = 0
$ws.Range("C24").Value2 = 186
$ws.Range("D24").Value2 = 158
$ws.Range("E24").Value2 = 17.721518987341
$ws.Range("F24").Value2 = 700
$ws.Range("G24").Value2 = 647
$ws.Range("H24").Value2 = 8.191653786707
$ws.Range("I24").Value2 = 509
$ws.Range("J24").Value2 = 481
$ws.Range("K24").Value2 = 5.821205821205
$ws.Range("L24").Value2 = 20.330969267139
$ws.Range("M24").Value2 = 39.071038251366
$ws.Range("C25").Value2 = 79
$ws.Range("D25").Value2 = 57
$ws.Range("E25").Value2 = 38.596491228070
$ws.Range("F25").Value2 = 324
$ws.Range("G25").Value2 = 231
$ws.Range("H25").Value2 = 40.259740259740
$ws.Range("I25").Value2 = 253
$ws.Range("J25").Value2 = 181
$ws.Range("K25").Value2 = 39.779005524861
$ws.Range("L25").Value2 = 58.125
$ws.Range("M25").Value2 = -4.887218045112
$ws.Range("C26").Value2 = 4
$ws.Range("D26").Value2 = 5
$ws.Range("E26").Value2 = -20
$ws.Range("F26").Value2 = 18
$ws.Range("G26").Value2 = 14
$ws.Range("H26").Value2 = 28.571428571428
$ws.Range("I26").Value2 = 14
$ws.Range("J26").Value2 = 14
$ws.Range("K26").Value2 = 0
$ws.Range("L26").Value2 = 27.272727272727
$ws.Range("C27").Value2 = 8
$ws.Range("D27").Value2 = 6
$ws.Range("E27").Value2 = 33.333333333333
$ws.Range("F27").Value2 = 35
$ws.Range("G27").Value2 = 21
$ws.Range("H27").Value2 = 66.666666666666
$ws.Range("I27").Value2 = 28
$ws.Range("J27").Value2 = 19
$ws.Range("K27").Value2 = 47.368421052631
$ws.Range("L27").Value2 = 27.272727272727
$ws.Range("C28").Value2 = 9
$ws.Range("E28").Value2 = 800
$ws.Range("F28").Value2 = 15
$ws.Range("G28").Value2 = 7
$ws.Range("H28").Value2 = 114.285714285714
$ws.Range("I28").Value2 = 12
$ws.Range("J28").Value2 = 3
$ws.Range("K28").Value2 = 300
$ws.Range("L28").Value2 = 50
$ws.Range("M28").Value2 = 200
$ws.Range("N28").Value2 = -70
$ws.Range("C29").Value2 = 5
$ws.Range("E29").Value2 = 400
$ws.Range("F29").Value2 = 10
$ws.Range("G29").Value2 = 6
$ws.Range("H29").Value2 = 66.666666666666
$ws.Range("I29").Value2 = 7
$ws.Range("J29").Value2 = 3
$ws.Range("K29").Value2 = 133.333333333333
$ws.Range("L29").Value2 = 40
$ws.Range("M29").Value2 = 133.333333333333
$ws.Range("N29").Value2 = -81.081081081081
$ws.Range("F30").Value2 = 3
$ws.Range("H30").Value2 = 50
$ws.Range("I30").Value2 = 3
$ws.Range("K30").Value2 = 200

# --- Row 30 ("Hate Crimes") column swap ------------------------------------
# C30: was the text placeholder "0" -> now a real count (switches to the
#      integer style used by sibling cells).
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value2 = 1

# D30: was a real count -> now the text placeholder "0" (switches to the
#      general/text style). Format as text first so the numeric-looking
#      "0" is stored as a string instead of being coerced to a number, then
#      restore the General number format (lower-case so the engine reuses
#      the existing style instead of minting a duplicate).
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0"
$ws.Range("D30").NumberFormat = "general"

# E30: was a real % change -> now the text placeholder "***.*".
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "***.*"
$ws.Range("E30").NumberFormat = "general"
